$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.1434
$ws.Range("A14").Value = -21.89279999999999
$ws.Range("A21").Value = -20.10949999999999
$ws.Range("D22").Value = -8.070200000000003
$ws.Range("A23").Value = -20.30899999999998
$ws.Range("D24").Value = -7.600399999999997
$ws.Range("A25").Value = -21.76209999999999
$ws.Range("A26").Value = -21.05569999999997
$ws.Range("D28").Value = -7.943099999999996
$ws.Range("A29").Value = -21.07369999999998
$ws.Range("D36").Value = -6.991799999999997
$ws.Range("D45").Value = -7.200999999999997
$ws.Range("D48").Value = -7.267699999999994
$ws.Range("D49").Value = -7.890799999999998
$ws.Range("D52").Value = -7.976000000000007
$ws.Range("A53").Value = -21.8973
$ws.Range("D53").Value = -7.924299999999999
$ws.Range("D54").Value = -8.065700000000003
$ws.Range("A57").Value = -22.30050000000001
$ws.Range("A59").Value = -21.92239999999999
$ws.Range("A69").Value = -21.61429999999998
$ws.Range("D70").Value = -6.581399999999999
$ws.Range("A79").Value = -20.3673
$ws.Range("A83").Value = -21.86650000000001
$ws.Range("D86").Value = -8.727399999999999
$ws.Range("D87").Value = -8.622299999999997
$ws.Range("D89").Value = -8.988899999999999
$ws.Range("A91").Value = -20.71819999999997
$ws.Range("A93").Value = -21.50210000000001
$ws.Range("D101").Value = -8.048200000000003
$ws.Range("A103").Value = -21.7132
